$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 414.25
$v2 = $ws.Cells.Item(2, 8).Value2
Write-Host $v2
$v3 = $ws.Range("H2").Value2
Write-Host $v3
